# Actualización automática 2025-06-02 14:06:09
# Adds a "PRESUPUESTO" (budget) column (G) to the "VENTA MENSUAL" sheet,
# filled with 0 for every data row and the totals row, copying the
# formatting from the neighboring "junio" column (F) so the new column
# reuses the existing cell styles instead of creating new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Header cell: reuse the bold/centered header style from F1
$ws.Range("G1").Value = "PRESUPUESTO"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Data rows 2-18: numeric 0, reuse the currency style from column F
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 6).Copy()
    $ws.Cells.Item($r, 7).PasteSpecial(-4122)
}

# Totals row 19: numeric 0, reuse the totals style from F19
$ws.Range("G19").Value = 0
$ws.Range("F19").Copy()
$ws.Range("G19").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Set the new column width to match the diff (stored XML width of 17).
# Excel's COM ColumnWidth is reported ~0.83 narrower than the raw OOXML
# <col width> value (font padding offset), so compensate here.
$ws.Columns.Item(7).ColumnWidth = 16.17
